$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 186.99
$ws.Range("I15").Value = 186.99
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 560.97
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -391.97

$ws.Range("H19").Value = 1524
$ws.Range("I19").Value = 1867
$ws.Range("J19").Value = 1377
$ws.Range("K19").Value = 1867
$ws.Range("L19").Value = 1377
$ws.Range("M19").Value = -1692
$ws.Range("N19").Value = -1727

$ws.Range("H28").Value = 1109.5385
$ws.Range("I28").Value = 324.8889
$ws.Range("J28").Value = 2875
$ws.Range("K28").Value = 324.8889
$ws.Range("L28").Value = 2875
$ws.Range("M28").Value = 160.1111
$ws.Range("N28").Value = -3845

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").ClearContents()

$ws.Range("H53").Value = 283.93103
$ws.Range("I53").Value = 227.53334
$ws.Range("J53").Value = 344.35715
$ws.Range("K53").Value = 227.53334
$ws.Range("L53").Value = 344.35715
$ws.Range("M53").Value = 409.46666
$ws.Range("N53").Value = -1618.35715

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()

$ws.Range("H129").Value = 2531.8813
$ws.Range("I129").Value = 11560.667
$ws.Range("J129").Value = 906.7
$ws.Range("K129").Value = 34682.001
$ws.Range("L129").Value = 2720.1
$ws.Range("M129").Value = -29682.001
$ws.Range("N129").Value = -12720.1

$ws.Range("H137").Value = 1398.878
$ws.Range("I137").Value = 1038.6857
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 3116.0571
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -566.0571
$ws.Range("N137").Value = -15600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 167865.17
$ws.Range("I2").Value = 1647.75
$ws.Range("J2").Value = 500300
$ws.Range("K2").Value = 1647.75
$ws.Range("L2").Value = 500300
$ws.Range("M2").Value = -1534.75
$ws.Range("N2").Value = -500526

$ws.Range("H32").Value = 27539.012
$ws.Range("I32").Value = 10320.822
$ws.Range("J32").Value = 111334.2
$ws.Range("K32").Value = 10320.822
$ws.Range("L32").Value = 111334.2
$ws.Range("M32").Value = -10033.822
$ws.Range("N32").Value = -111908.2

$ws.Range("H61").Value = 2050.3333
$ws.Range("I61").Value = 1720
$ws.Range("J61").Value = 2103.6128
$ws.Range("K61").Value = 1720
$ws.Range("L61").Value = 2103.6128
$ws.Range("M61").Value = -1508
$ws.Range("N61").Value = -2527.6128

$ws.Range("H110").Value = 77078400
$ws.Range("I110").Value = 77078400
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 77078400
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -77076355
$ws.Range("N110").ClearContents()

$ws.Range("H116").Value = 167865.17
$ws.Range("I116").Value = 1647.75
$ws.Range("J116").Value = 500300
$ws.Range("K116").Value = 1647.75
$ws.Range("L116").Value = 500300
$ws.Range("M116").Value = 646.25
$ws.Range("N116").Value = -504888

$ws.Range("H132").Value = 16376.333
$ws.Range("I132").Value = 19441.412
$ws.Range("J132").Value = 3349.75
$ws.Range("K132").Value = 58324.236
$ws.Range("L132").Value = 10049.25
$ws.Range("M132").Value = -55794.236

$ws.Range("H136").Value = 2050.3333
$ws.Range("I136").Value = 1720
$ws.Range("J136").Value = 2103.6128
$ws.Range("K136").Value = 5160
$ws.Range("L136").Value = 6310.8384
$ws.Range("M136").Value = -2610
$ws.Range("N136").Value = -11410.8384

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 167865.17
$ws.Range("I3").Value = 1647.75
$ws.Range("J3").Value = 500300
$ws.Range("K3").Value = 1647.75
$ws.Range("L3").Value = 500300
$ws.Range("M3").Value = -1533.75
$ws.Range("N3").Value = -500528

$ws.Range("H134").Value = 2855.7322
$ws.Range("I134").Value = 2762.6365
$ws.Range("J134").Value = 3197.0833
$ws.Range("K134").Value = 8287.9095
$ws.Range("L134").Value = 9591.249899999999
$ws.Range("M134").Value = -5752.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 173.13333
$ws.Range("I7").Value = 89.5
$ws.Range("J7").Value = 340.4
$ws.Range("K7").Value = 89.5
$ws.Range("L7").Value = 340.4
$ws.Range("M7").Value = 23.5
$ws.Range("N7").Value = -566.4

$ws.Range("H10").Value = 400
$ws.Range("I10").Value = 400
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 400
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -261

$ws.Range("H31").Value = 31040.44
$ws.Range("I31").Value = 1070.9474
$ws.Range("J31").Value = 49408.84
$ws.Range("K31").Value = 1070.9474
$ws.Range("L31").Value = 49408.84
$ws.Range("M31").Value = -775.9474
$ws.Range("N31").Value = -49998.84

$ws.Range("H34").Value = 31040.44
$ws.Range("I34").Value = 1070.9474
$ws.Range("J34").Value = 49408.84
$ws.Range("K34").Value = 1070.9474
$ws.Range("L34").Value = 49408.84
$ws.Range("M34").Value = -868.9474
$ws.Range("N34").Value = -49812.84

$ws.Range("H58").Value = 2243.8
$ws.Range("I58").Value = 2134.4285
$ws.Range("J58").Value = 2499
$ws.Range("K58").Value = 2134.4285
$ws.Range("L58").Value = 2499
$ws.Range("M58").Value = -1931.4285
$ws.Range("N58").Value = -2905

$ws.Range("H96").Value = 23500
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 23500
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 23500
$ws.Range("N96").Value = -28992

$ws.Range("H99").Value = 2245
$ws.Range("I99").Value = 1898.7142
$ws.Range("J99").Value = 2387.5881
$ws.Range("K99").Value = 1898.7142
$ws.Range("L99").Value = 2387.5881
$ws.Range("M99").Value = -400.7141999999999

$ws.Range("H126").Value = 2245
$ws.Range("I126").Value = 1898.7142
$ws.Range("J126").Value = 2387.5881
$ws.Range("K126").Value = 5696.142599999999
$ws.Range("L126").Value = 7162.7643
$ws.Range("M126").Value = -3226.142599999999

$ws.Range("H136").Value = 2243.8
$ws.Range("I136").Value = 2134.4285
$ws.Range("J136").Value = 2499
$ws.Range("K136").Value = 6403.2855
$ws.Range("L136").Value = 7497
$ws.Range("M136").Value = -3853.2855
$ws.Range("N136").Value = -12597

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 120.73333
$ws.Range("I8").Value = 120.73333
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 362.19999
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -223.19999

$ws.Range("H23").Value = 477.70587
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 477.70587
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1433.11761
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1903.11761

$ws.Range("H38").Value = 63.76923
$ws.Range("I38").Value = 68
$ws.Range("J38").Value = 61.125
$ws.Range("K38").Value = 204
$ws.Range("L38").Value = 183.375
$ws.Range("M38").Value = 143
$ws.Range("N38").Value = -877.375

$ws.Range("H56").Value = 4249.4165
$ws.Range("I56").Value = 4249.4165
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4249.4165
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -3719.4165

$ws.Range("H131").Value = 799.7895
$ws.Range("I131").Value = 482.83334
$ws.Range("J131").Value = 873.8831
$ws.Range("K131").Value = 1448.50002
$ws.Range("L131").Value = 2621.6493
$ws.Range("M131").Value = 3591.49998
$ws.Range("N131").Value = -12701.6493

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 900.6923
$ws.Range("I22").Value = 1041.3334
$ws.Range("J22").Value = 708.9091
$ws.Range("K22").Value = 1041.3334
$ws.Range("L22").Value = 708.9091
$ws.Range("M22").Value = -746.3334
$ws.Range("N22").Value = -1298.9091

$ws.Range("H27").Value = 900.6923
$ws.Range("I27").Value = 1041.3334
$ws.Range("J27").Value = 708.9091
$ws.Range("K27").Value = 1041.3334
$ws.Range("L27").Value = 708.9091
$ws.Range("M27").Value = -934.3334
$ws.Range("N27").Value = -922.9091

$ws.Range("H46").Value = 7024.875
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 7742.7144
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 7742.7144
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -8118.7144

$ws.Range("H61").Value = 3421
$ws.Range("I61").Value = 3421
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3421
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3219
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 3421
$ws.Range("I113").Value = 3421
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3421
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1251
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 4273.522
$ws.Range("I132").Value = 6280.8
$ws.Range("J132").Value = 2729.4614
$ws.Range("K132").Value = 18842.4
$ws.Range("L132").Value = 8188.3842
$ws.Range("M132").Value = -16312.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3024.56
$ws.Range("I132").Value = 1537
$ws.Range("J132").Value = 7735.1665
$ws.Range("K132").Value = 4611
$ws.Range("L132").Value = 23205.4995
$ws.Range("M132").Value = -2081
$ws.Range("N132").Value = -28265.4995
